$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the `=TRUE()` boolean formulas in E2/E3 with the literal text
# "TRUE". Assigning the bare word "TRUE" via .Value gets auto-coerced into
# a real Boolean by Excel, so instead write it with a trailing space (which
# keeps it text), then use a scratch helper cell + TRIM()/paste-values to
# strip the space while staying on the text code path (instead of
# re-parsing the literal and re-coercing it back to a Boolean).
$ws.Range("E2").Value = "TRUE "
$ws.Range("E3").Value = "TRUE "

$helper = $ws.Range("G2:G3")
$helper.Formula = "=TRIM(E2)"
$helper.Copy()
$ws.Range("E2:E3").PasteSpecial(-4163)
$helper.Clear()

# Update selection to match the diff (E2:E3 active range)
$ws.Range("E2:E3").Select()
